# Powerpoint writer: consolidate text run nodes.
# Merge adjacent a:r runs that got split purely because of word-boundary
# spaces, without altering the rendered text. We do this by rewriting a
# character range that spans the boundary between two runs with the same
# text it already had -- the host recomputes the minimal set of runs
# needed to express that text, collapsing the old (now redundant) split.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: "Testing" + " " + "custom" + " " + "properties" ---
# becomes "Testing " + "custom " + "properties"
$title = $s.Shapes.Item(1).TextFrame.TextRange

# Merge "Testing" (chars 1-7) with the following space (char 8) -> "Testing "
$title.Characters(1, 8).Text = "Testing "

# Merge "custom" (chars 9-14) with the following space (char 15) -> "custom "
$title.Characters(9, 7).Text = "custom "

# --- Subtitle shape: two line breaks, then "A." + " " + "M." ---
# becomes two line breaks, then "A. " + "M."
$subtitle = $s.Shapes.Item(2).TextFrame.TextRange

# Characters 1-2 are the two leading <a:br/> breaks; "A." starts at 3.
# Merge "A." (chars 3-4) with the following space (char 5) -> "A. "
$subtitle.Characters(3, 3).Text = "A. "
